$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column D
$ws.Range("D1").Value = "ITI"

# Update ConditionType (column C) and add ITI (column D) values for rows 2-17
$data = @(
    @(4,6),
    @(4,8),
    @(2,6),
    @(2,6),
    @(3,7),
    @(3,6),
    @(2,8),
    @(1,8),
    @(1,8),
    @(2,7),
    @(4,8),
    @(3,6),
    @(3,8),
    @(1,7),
    @(4,6),
    @(1,7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $data[$i][0]
    $ws.Cells.Item($row, 4).Value = $data[$i][1]
}

# Remove the trailing rows (trials 17-19), which were rows 18-20
$ws.Range("A18:D20").Delete() | Out-Null

# Update selection to match target workbook state
$ws.Range("F5").Select() | Out-Null
